# Add data for 2024-06-08
# Updates the 2024 (column K) violent-crime counts across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and each affected individual neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 3304
$ws.Range('K3').Value = 3277
$ws.Range('K4').Value = 679
$ws.Range('K5').Value = 213
$ws.Range('K6').Value = 3864
$ws.Range('K7').Value = 11337

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 218
$ws.Range('K3').Value = 228
$ws.Range('K4').Value = 41
$ws.Range('K5').Value = 21
$ws.Range('K6').Value = 243
$ws.Range('K7').Value = 751

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 79
$ws.Range('K7').Value = 246

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 126
$ws.Range('K3').Value = 163
$ws.Range('K6').Value = 129
$ws.Range('K7').Value = 449

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K6').Value = 47
$ws.Range('K7').Value = 184

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 107
$ws.Range('K6').Value = 119
$ws.Range('K7').Value = 392

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 69
$ws.Range('K6').Value = 104
$ws.Range('K7').Value = 270

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K4').Value = 13
$ws.Range('K7').Value = 197

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 91
$ws.Range('K7').Value = 319
$ws.Range('K8').Value = 751
$ws.Range('K9').Value = 47
$ws.Range('K11').Value = 240
$ws.Range('K12').Value = 18
$ws.Range('K15').Value = 116
$ws.Range('K16').Value = 38
$ws.Range('K19').Value = 347
$ws.Range('K20').Value = 260
$ws.Range('K22').Value = 36
$ws.Range('K23').Value = 109
$ws.Range('K27').Value = 118
$ws.Range('K29').Value = 591
$ws.Range('K33').Value = 449
$ws.Range('K36').Value = 134
$ws.Range('K37').Value = 392
$ws.Range('K42').Value = 404
$ws.Range('K47').Value = 60
$ws.Range('K48').Value = 142
$ws.Range('K51').Value = 132
$ws.Range('K52').Value = 306
$ws.Range('K54').Value = 221
$ws.Range('K58').Value = 3
$ws.Range('K60').Value = 70
$ws.Range('K63').Value = 38
$ws.Range('K64').Value = 69
$ws.Range('K65').Value = 270
$ws.Range('K66').Value = 42
$ws.Range('K70').Value = 19
$ws.Range('K72').Value = 54
$ws.Range('K77').Value = 80
$ws.Range('K78').Value = 145
$ws.Range('K79').Value = 295
$ws.Range('K80').Value = 41
$ws.Range('K82').Value = 14
$ws.Range('K83').Value = 246
$ws.Range('K85').Value = 539
$ws.Range('K86').Value = 75
$ws.Range('K89').Value = 150
$ws.Range('K95').Value = 184
$ws.Range('K99').Value = 197
$ws.Range('K101').Value = 11337

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 65
$ws.Range('K6').Value = 103
$ws.Range('K7').Value = 221

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K4').Value = 34
$ws.Range('K7').Value = 591

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K2').Value = 21
$ws.Range('K6').Value = 76
$ws.Range('K7').Value = 142

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 117
$ws.Range('K3').Value = 91
$ws.Range('K6').Value = 112
$ws.Range('K7').Value = 347

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 104
$ws.Range('K3').Value = 129
$ws.Range('K7').Value = 404

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 43
$ws.Range('K7').Value = 145

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 37
$ws.Range('K3').Value = 36
$ws.Range('K6').Value = 28
$ws.Range('K7').Value = 109

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K4').Value = 16
$ws.Range('K7').Value = 295

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 26
$ws.Range('K7').Value = 69

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 73
$ws.Range('K7').Value = 260

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K3').Value = 40
$ws.Range('K7').Value = 134

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 112
$ws.Range('K7').Value = 319

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K3').Value = 26
$ws.Range('K6').Value = 60

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K6').Value = 18
$ws.Range('K7').Value = 60

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K3').Value = 29
$ws.Range('K7').Value = 116

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('K6').Value = 23
$ws.Range('K7').Value = 42

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K6').Value = 93
$ws.Range('K7').Value = 240

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K6').Value = 10
$ws.Range('K7').Value = 47

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 25
$ws.Range('K4').Value = 8
$ws.Range('K6').Value = 31
$ws.Range('K7').Value = 91

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('K2').Value = 7
$ws.Range('K7').Value = 19

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 36
$ws.Range('K7').Value = 150

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K4').Value = 13
$ws.Range('K7').Value = 118

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 28
$ws.Range('K7').Value = 75

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 38
$ws.Range('K7').Value = 132

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K2').Value = 20
$ws.Range('K6').Value = 20
$ws.Range('K7').Value = 70

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 191
$ws.Range('K3').Value = 184
$ws.Range('K4').Value = 28
$ws.Range('K6').Value = 124
$ws.Range('K7').Value = 539

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('K2').Value = 19
$ws.Range('K7').Value = 36

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K3').Value = 17
$ws.Range('K7').Value = 54

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('K3').Value = 4
$ws.Range('K6').Value = 14

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K3').Value = 29
$ws.Range('K7').Value = 80

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K2').Value = 11
$ws.Range('K7').Value = 41

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 81
$ws.Range('K6').Value = 122
$ws.Range('K7').Value = 306

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('K6').Value = 4
$ws.Range('K7').Value = 18

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K2').Value = 10
$ws.Range('K7').Value = 38

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range('K6').Value = 3
$ws.Range('K7').Value = 3
